$wb = $excel.ActiveWorkbook

# --- Metrics sheet: update raw metric values (B2:B13) ---
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 220715.43
$metrics.Range("B3").Value  = 181186.62000000002
$metrics.Range("B4").Value  = 70183.360000000001
$metrics.Range("B5").Value  = 8760
$metrics.Range("B6").Value  = 4587846.8999999994
$metrics.Range("B7").Value  = 3871005.29
$metrics.Range("B8").Value  = 1340785.5000000002
$metrics.Range("B9").Value  = 177761
$metrics.Range("B10").Value = 33053170.70099983
$metrics.Range("B11").Value = 19900875.360000003
$metrics.Range("B12").Value = 11622494.390000002
$metrics.Range("B13").Value = 1275388

# Move the saved selection on the Metrics sheet to D8 (as recorded in the workbook view)
$metrics.Range("D8").Select()

# --- today sheet: move the saved selection to I11 ---
$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("I11").Select()
